$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 40900
$ws.Range("J3").Value = 40900
$ws.Range("L3").Value = 40900
$ws.Range("N3").Value = -41128
$ws.Range("H6").Value = 119
$ws.Range("I6").Value = 108.84615
$ws.Range("K6").Value = 326.53845
$ws.Range("M6").Value = -214.53845
$ws.Range("H41").Value = 2340
$ws.Range("I41").Value = 466.66666
$ws.Range("K41").Value = 466.66666
$ws.Range("M41").Value = -26.66665999999998
$ws.Range("H96").Value = 482.30768
$ws.Range("I96").Value = 88.71429000000001
$ws.Range("J96").Value = 941.5
$ws.Range("K96").Value = 266.14287
$ws.Range("L96").Value = 2824.5
$ws.Range("M96").Value = 1106.85713
$ws.Range("N96").Value = -5570.5
$ws.Range("H102").Value = 40900
$ws.Range("J102").Value = 40900
$ws.Range("L102").Value = 40900
$ws.Range("N102").Value = -47390
$ws.Range("H125").Value = 1882.1111
$ws.Range("I125").Value = 1862.8572
$ws.Range("K125").Value = 16765.7148
$ws.Range("M125").Value = -14305.7148
$ws.Range("H127").Value = 1741.3334
$ws.Range("I127").Value = 834
$ws.Range("K127").Value = 2502
$ws.Range("M127").Value = 2458
$ws.Range("H129").Value = 3664.6667
$ws.Range("J129").Value = 9000
$ws.Range("L129").Value = 27000
$ws.Range("N129").Value = -37000
$ws.Range("H137").Value = 3245.8823
$ws.Range("I137").Value = 2147.75
$ws.Range("K137").Value = 6443.25
$ws.Range("M137").Value = -3893.25
$ws.Range("H138").Value = 2429.1428
$ws.Range("J138").Value = 2662.5
$ws.Range("L138").Value = 7987.5
$ws.Range("N138").Value = -18267.5

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2382630.8
$ws.Range("I32").Value = 499.7143
$ws.Range("K32").Value = 499.7143
$ws.Range("M32").Value = -212.7143
$ws.Range("H45").Value = 2677.0667
$ws.Range("I45").Value = 2192.4
$ws.Range("K45").Value = 2192.4
$ws.Range("M45").Value = -1815.4
$ws.Range("H61").Value = 3849.0625
$ws.Range("I61").Value = 2509.7778
$ws.Range("J61").Value = 5571
$ws.Range("K61").Value = 2509.7778
$ws.Range("L61").Value = 5571
$ws.Range("M61").Value = -2297.7778
$ws.Range("N61").Value = -5995
$ws.Range("H122").Value = 2961.5833
$ws.Range("I122").Value = 2845.3635
$ws.Range("K122").Value = 8536.0905
$ws.Range("M122").Value = -6086.0905
$ws.Range("H136").Value = 3849.0625
$ws.Range("I136").Value = 2509.7778
$ws.Range("J136").Value = 5571
$ws.Range("K136").Value = 7529.3334
$ws.Range("L136").Value = 16713
$ws.Range("M136").Value = -4979.3334
$ws.Range("N136").Value = -21813

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H37").Value = 1624.625
$ws.Range("I37").Value = 1859.6
$ws.Range("J37").Value = 1233
$ws.Range("K37").Value = 1859.6
$ws.Range("L37").Value = 1233
$ws.Range("M37").Value = -1722.6
$ws.Range("N37").Value = -1507
$ws.Range("H94").Value = 1476.0454
$ws.Range("I94").Value = 1485.2106
$ws.Range("J94").Value = 1418
$ws.Range("K94").Value = 1485.2106
$ws.Range("L94").Value = 1418
$ws.Range("M94").Value = -1034.2106
$ws.Range("N94").Value = -2320
$ws.Range("H105").Value = 2684
$ws.Range("I105").Value = 2641.3333
$ws.Range("K105").Value = 2641.3333
$ws.Range("M105").Value = -894.3332999999998
$ws.Range("H134").Value = 1423.0625
$ws.Range("I134").Value = 1423.0625
$ws.Range("K134").Value = 4269.1875
$ws.Range("M134").Value = -1734.1875

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 57600
$ws.Range("J4").Value = 505000
$ws.Range("L4").Value = 505000
$ws.Range("N4").Value = -505224
$ws.Range("H7").Value = 3084.7837
$ws.Range("I7").Value = 5726.9443
$ws.Range("J7").Value = 581.6842
$ws.Range("K7").Value = 5726.9443
$ws.Range("L7").Value = 581.6842
$ws.Range("M7").Value = -5613.9443
$ws.Range("N7").Value = -807.6842
$ws.Range("H31").Value = 6436
$ws.Range("I31").Value = 1058.625
$ws.Range("K31").Value = 1058.625
$ws.Range("M31").Value = -763.625
$ws.Range("H34").Value = 6436
$ws.Range("I34").Value = 1058.625
$ws.Range("K34").Value = 1058.625
$ws.Range("M34").Value = -856.625
$ws.Range("H58").Value = 5004.2
$ws.Range("J58").Value = 5336.3335
$ws.Range("L58").Value = 5336.3335
$ws.Range("N58").Value = -5742.3335
$ws.Range("H132").Value = 3560.2856
$ws.Range("I132").Value = 3167.6365
$ws.Range("K132").Value = 9502.9095
$ws.Range("M132").Value = -6972.9095
$ws.Range("H136").Value = 5004.2
$ws.Range("J136").Value = 5336.3335
$ws.Range("L136").Value = 16009.0005
$ws.Range("N136").Value = -21109.0005

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1011.3333
$ws.Range("J34").Value = 1406.125
$ws.Range("L34").Value = 4218.375
$ws.Range("N34").Value = -4386.375
$ws.Range("H132").Value = 1858.6471
$ws.Range("I132").Value = 1877.4445
$ws.Range("J132").Value = 1837.5
$ws.Range("K132").Value = 16897.0005
$ws.Range("L132").Value = 16537.5
$ws.Range("M132").Value = -14367.0005
$ws.Range("N132").Value = -21597.5

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1465.7858
$ws.Range("I102").Value = 1556.5454
$ws.Range("J102").Value = 1133
$ws.Range("K102").Value = 1556.5454
$ws.Range("L102").Value = 1133
$ws.Range("M102").Value = 65.45460000000003
$ws.Range("N102").Value = -4377
$ws.Range("H132").Value = 3935.875
$ws.Range("I132").Value = 2988.5557
$ws.Range("K132").Value = 8965.667099999999
$ws.Range("M132").Value = -6435.667099999999

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2813.0715
$ws.Range("I61").Value = 1398.5454
$ws.Range("K61").Value = 1398.5454
$ws.Range("M61").Value = -1196.5454
$ws.Range("H98").Value = 65000
$ws.Range("J98").Value = 65000
$ws.Range("L98").Value = 65000
$ws.Range("N98").Value = -70990
$ws.Range("H100").Value = 9999.857
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 9999.857
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 9999.857
$ws.Range("N100").Value = -11081.857
$ws.Range("M100").ClearContents()
$ws.Range("H113").Value = 2813.0715
$ws.Range("I113").Value = 1398.5454
$ws.Range("K113").Value = 1398.5454
$ws.Range("M113").Value = 771.4546
$ws.Range("H132").Value = 4200
$ws.Range("I132").Value = 3400
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 10200
$ws.Range("L132").Value = 15000
$ws.Range("M132").Value = -7670
$ws.Range("N132").Value = -20060

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H28").Value = 20000
$ws.Range("J28").Value = 20000
$ws.Range("L28").Value = 20000
$ws.Range("N28").Value = -20696
$ws.Range("H31").Value = 14333.333
$ws.Range("J31").Value = 14333.333
$ws.Range("L31").Value = 14333.333
$ws.Range("N31").Value = -15029.333
$ws.Range("H100").Value = 1215.7273
$ws.Range("I100").Value = 1137.3
$ws.Range("K100").Value = 2274.6
$ws.Range("M100").Value = -1733.6
